$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data occupies rows 2..49 (4 years x 12 months), columns A..F.
# Each 12-row block (one year) needs its rows rotated so that the last
# three months (Oct, Nov, Dec) move to the front of the block, followed
# by Jan..Sep, e.g. for 2014 (rows 2-13):
#   before: 2014-01 .. 2014-12
#   after : 2014-10, 2014-11, 2014-12, 2014-01, .., 2014-09

$firstRow = 2
$blockSize = 12
$lastRow = 49
$numBlocks = [int](($lastRow - $firstRow + 1) / $blockSize)
$cols = @("A", "B", "C", "D", "E", "F")

for ($b = 0; $b -lt $numBlocks; $b++) {
    $blockStart = $firstRow + $b * $blockSize

    # Snapshot the 12 rows of this block (as values) before rewriting.
    $snapshot = @()
    for ($i = 0; $i -lt $blockSize; $i++) {
        $r = $blockStart + $i
        $rowVals = @()
        foreach ($col in $cols) {
            $rowVals += $ws.Range($col + $r).Value()
        }
        $snapshot += ,$rowVals
    }

    # New row order within the block: indices 9,10,11 (Oct,Nov,Dec) first,
    # then indices 0..8 (Jan..Sep).
    $newOrder = @(9, 10, 11, 0, 1, 2, 3, 4, 5, 6, 7, 8)

    for ($i = 0; $i -lt $blockSize; $i++) {
        $r = $blockStart + $i
        $srcVals = $snapshot[$newOrder[$i]]
        for ($c = 0; $c -lt $cols.Length; $c++) {
            $ws.Range($cols[$c] + $r).Value = $srcVals[$c]
        }
    }
}
